$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8530270457267761
$ws.Range("B1").Value = 1.274253129959106
$ws.Range("C1").Value = 2.668391704559326
$ws.Range("D1").Value = 3.936103343963623
$ws.Range("E1").Value = 1.849095821380615
